$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new category/feature values in column D for rows 10-11,
# mirroring the existing "Строительство" block structure (rows 10-12)
# the same way column D mirrors column C for the "Население" block
# (rows 3-5).

# Row 10 (header row): same formatting as C10 (style used for headers)
$ws.Range("C10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").Value = "Уровень жизни"

# Row 11 (data row): same formatting as D5 (the matching style used for
# the single-value rows under column D)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = "Жил. площ.на одного чел. - livarea (кв. м) (id8211001)"

# Row 12 stays empty but keeps the same cell formatting as D11/D5
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the active selection from D20 to D14
$ws.Range("D14").Select() | Out-Null
